# Added covariance for some assets:
# Remove the VTI and ITOT rows from the "Data" sheet, shifting the
# remaining rows (S&P 500, SAM 100 ... SAM 30/70) up by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Activate() | Out-Null

# Delete entire rows 2 (VTI) and 3 (ITOT); deleting row 2 twice shifts
# row 3's original content into row 2, so deleting row 2 twice removes
# both original rows.
$ws.Rows.Item(2).Delete() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null

# Update the active selection to match the post-edit state (row 2 is
# selected after the deletion).
$ws.Range("A2:XFD2").Select() | Out-Null
